# edit.ps1
# Applies the "updated names on slides." commit:
#   - Slide 10 ("Thank you !" title): split the single run into two runs
#       " Thank " + "you !" (same combined text, but a leading space is
#       introduced before "Thank").
#   - Slide 6 Rectangle-3 name tag: "Vijay Kumar Karanam"
#       -> "Prajakt Uttamrao Khawase" (3 runs) + reposition/resize box.
#   - Slide 7 Rectangle-3 name tag: "Vijay Kumar Karanam"
#       -> "Harsha Vardhan Reddy Nallavolu" (2 runs) + reposition/resize box.
#   - Slide 8 Rectangle-3 name tag: "Vijay Kumar Karanam"
#       -> "Harsha Vardhan Reddy Nallavolu" (2 runs) + reposition/resize box.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 - "Thank you !" title -> " Thank " / "you !" (two runs)
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$title = $s10.Shapes.Item(1)
$tr10 = $title.TextFrame.TextRange

# Original text is "Thank you !" (11 chars). Re-point the first 6
# characters ("Thank ") to " Thank " - this both prepends the leading
# space and forces PowerPoint to split the text into two runs, matching
# the target XML (first run " Thank ", second run "you !").
$tr10.Characters(1, 6).Text = " Thank "

# ---------------------------------------------------------------------------
# Slide 6 - Rectangle 3 name tag
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$rect6 = $s6.Shapes.Item(3)

$tr6 = $rect6.TextFrame.TextRange
$tr6.Text = ""
$tr6.InsertAfter("Prajakt") | Out-Null
$tr6.InsertAfter(" Uttamrao ") | Out-Null
$tr6.InsertAfter("Khawase") | Out-Null

$rect6.Left = 707.875
$rect6.Top = 478.0
$rect6.Width = 222.49685668945312

# ---------------------------------------------------------------------------
# Slide 7 - Rectangle 3 name tag
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$rect7 = $s7.Shapes.Item(3)

$tr7 = $rect7.TextFrame.TextRange
$tr7.Text = ""
$tr7.InsertAfter("Harsha") | Out-Null
$tr7.InsertAfter(" Vardhan Reddy Nallavolu") | Out-Null

$rect7.Left = 659.875
$rect7.Top = 463.4593811035156
$rect7.Width = 272.6569519042969

# ---------------------------------------------------------------------------
# Slide 8 - Rectangle 3 name tag
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$rect8 = $s8.Shapes.Item(3)

$tr8 = $rect8.TextFrame.TextRange
$tr8.Text = ""
$tr8.InsertAfter("Harsha") | Out-Null
$tr8.InsertAfter(" Vardhan Reddy Nallavolu") | Out-Null

$rect8.Left = 659.875
$rect8.Top = 466.9759216308594
$rect8.Width = 272.6569519042969
